# Apply the changes described by the commit: refresh "want-to-go" counters
# on sheet "展览" and "全部类型", and insert a newly scraped event as the
# new row 41 on sheet "展览" (pushing the existing rows 41-43 down to 42-44).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# ---------------------------------------------------------------------
# 1) Sheet "展览" - update the "想去人数" (F column) counters that changed
#    for rows that are not affected by the later row insertion (rows 1-40).
# ---------------------------------------------------------------------
$ws1.Range("F2").Value = 2614
$ws1.Range("F4").Value = 349
$ws1.Range("F5").Value = 1460
$ws1.Range("F6").Value = 1129
$ws1.Range("F13").Value = 8962
$ws1.Range("F15").Value = 2501
$ws1.Range("F19").Value = 615
$ws1.Range("F23").Value = 2065
$ws1.Range("F24").Value = 2148
$ws1.Range("F26").Value = 1851
$ws1.Range("F29").Value = 481
$ws1.Range("F30").Value = 822
$ws1.Range("F31").Value = 63
$ws1.Range("F32").Value = 138
$ws1.Range("F37").Value = 279
$ws1.Range("F38").Value = 463
$ws1.Range("F39").Value = 1332
$ws1.Range("F40").Value = 280

# ---------------------------------------------------------------------
# 2) Sheet "展览" - insert a brand new event as row 41, shifting the
#    previously-last three rows (41-43) down to 42-44.
# ---------------------------------------------------------------------
$ws1.Rows(41).Insert()

# restore the bordered/bold/centered number style on the new A41 cell by
# copying the formatting from the row right below it (which used to be
# row 41 before the shift, and still carries the correct style).
$ws1.Range("A42").Copy()
$ws1.Range("A41").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("A41").Value = 40
$ws1.Range("B41").NumberFormat = "@"
$ws1.Range("B41").Value = "2024-06-08"
$ws1.Range("C41").Value = "杭州·第八届YH樱花动漫游戏文化节"
$ws1.Range("D41").Value = "德胜东路2539号 梦马汽车小镇"
$ws1.Range("E41").Value = "2024.06.08 10:00-06.10 17:00"
$ws1.Range("F41").Value = 1
$ws1.Range("G41").Value = "不可售"
$ws1.Range("H41").Value = "https://show.bilibili.com/platform/detail.html?id=82687"
$ws1.Range("I41").Value = "//i2.hdslb.com/bfs/openplatform/202403/S5pnadXj1710210939138.png"

# Rows() Insert() only moves cells, it does not renumber the plain values
# that used to be in column A, so the running index in column A has to be
# bumped by one for every row that got pushed down.
$ws1.Range("A42").Value = 41
$ws1.Range("A43").Value = 42
$ws1.Range("A44").Value = 43

# the event that used to live on row 41 now lives on row 42, with an
# updated "want-to-go" counter
$ws1.Range("F42").Value = 74

# the event that used to live on row 43 now lives on row 44, with an
# updated "want-to-go" counter
$ws1.Range("F44").Value = 284

# ---------------------------------------------------------------------
# 3) Sheet "全部类型" - update the "想去人数" (F column) counters.
# ---------------------------------------------------------------------
$ws4.Range("F2").Value = 2614
$ws4.Range("F4").Value = 349
$ws4.Range("F5").Value = 1460
$ws4.Range("F7").Value = 1129
$ws4.Range("F13").Value = 8962
$ws4.Range("F15").Value = 2501
$ws4.Range("F20").Value = 615
$ws4.Range("F24").Value = 2065
$ws4.Range("F25").Value = 2148
$ws4.Range("F27").Value = 1851
$ws4.Range("F30").Value = 481
$ws4.Range("F31").Value = 822
$ws4.Range("F32").Value = 63
$ws4.Range("F33").Value = 138
$ws4.Range("F38").Value = 279
$ws4.Range("F39").Value = 463
$ws4.Range("F44").Value = 1332
$ws4.Range("F46").Value = 280
$ws4.Range("F47").Value = 74
$ws4.Range("F49").Value = 284
